$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A1").Value = "country"
$ws.Range("B1").Value = "channel"
$ws.Range("D1").Value = "product_link"

$null = $ws.Range("F18").Select()
